$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-23 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-24 Sunday", 2)
$d.Content.Find.Execute("89÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷6=", 2)
$d.Content.Find.Execute("83÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷6=", 2)
$d.Content.Find.Execute("47÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷4=", 2)
$d.Content.Find.Execute("35÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷6=", 2)
$d.Content.Find.Execute("15÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷7=", 2)
$d.Content.Find.Execute("76÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷6=", 2)
$d.Content.Find.Execute("54÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷2=", 2)
$d.Content.Find.Execute("88÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷5=", 2)
$d.Content.Find.Execute("31÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "77÷4=", 2)
$d.Content.Find.Execute("33÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷5=", 2)
$d.Content.Find.Execute("78÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷2=", 2)
$d.Content.Find.Execute("59÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷5=", 2)
$d.Content.Find.Execute("53÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷6=", 2)
$d.Content.Find.Execute("72÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷9=", 2)
$d.Content.Find.Execute("38÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷5=", 2)
$d.Content.Find.Execute("18÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "91÷5=", 2)
$d.Content.Find.Execute("42÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷2=", 2)
$d.Content.Find.Execute("32÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷9=", 2)
$d.Content.Find.Execute("87÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷4=", 2)
$d.Content.Find.Execute("81÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷4=", 2)
$d.Content.Find.Execute("84÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷5=", 2)
$d.Content.Find.Execute("20÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷8=", 2)
$d.Content.Find.Execute("31÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷3=", 2)
$d.Content.Find.Execute("10÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷7=", 2)
$d.Content.Find.Execute("68÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷7=", 2)
